$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 5

$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 150

$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 7
$ws.Range("G6").Value = 8
